$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (A1:G1) ---
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- Update D/E/F data values (rows 2-37) ---
$ws.Range("D2").Value = [double]"0.05517737733333333"
$ws.Range("E2").Value = [double]"0.52910956"
$ws.Range("F2").Value = [double]"1.5384917e-06"
$ws.Range("D3").Value = [double]"0.1791537466666667"
$ws.Range("E3").Value = [double]"3.2512784"
$ws.Range("F3").Value = [double]"4.9952819e-06"
$ws.Range("D4").Value = [double]"954.6141333333333"
$ws.Range("E4").Value = [double]"14031.106"
$ws.Range("F4").Value = [double]"0.026617176"
$ws.Range("D5").Value = [double]"438.1869266666666"
$ws.Range("E5").Value = [double]"6379.2934"
$ws.Range("F5").Value = [double]"0.012217814"
$ws.Range("D6").Value = [double]"0.15117088"
$ws.Range("E6").Value = [double]"3.7461084"
$ws.Range("F6").Value = [double]"4.2150453e-06"
$ws.Range("D7").Value = [double]"3.184314533333334"
$ws.Range("E7").Value = [double]"37.029248"
$ws.Range("F7").Value = [double]"8.8787139e-05"
$ws.Range("D8").Value = [double]"326.4253266666667"
$ws.Range("E8").Value = [double]"4821.4268"
$ws.Range("F8").Value = [double]"0.0091016043"
$ws.Range("D9").Value = [double]"0.8462718666666668"
$ws.Range("E9").Value = [double]"15.776121"
$ws.Range("F9").Value = [double]"2.3596305e-05"
$ws.Range("D10").Value = [double]"1.953059066666667"
$ws.Range("E10").Value = [double]"18.430519"
$ws.Range("F10").Value = [double]"5.4456469e-05"
$ws.Range("D11").Value = [double]"85.22702666666666"
$ws.Range("E11").Value = [double]"1091.9692"
$ws.Range("F11").Value = [double]"0.0023763556"
$ws.Range("D12").Value = [double]"8.864813999999999"
$ws.Range("E12").Value = [double]"117.56442"
$ws.Range("F12").Value = [double]"0.00024717453"
$ws.Range("D13").Value = [double]"15.50034"
$ws.Range("E13").Value = [double]"185.57011"
$ws.Range("F13").Value = [double]"0.0004321906"
$ws.Range("D14").Value = [double]"2.688615333333333"
$ws.Range("E14").Value = [double]"39.836747"
$ws.Range("F14").Value = [double]"7.496572999999999e-05"
$ws.Range("D15").Value = [double]"0.213"
$ws.Range("E15").Value = [double]"3.462316"
$ws.Range("F15").Value = [double]"5.9390052e-06"
$ws.Range("D16").Value = [double]"0.41"
$ws.Range("E16").Value = [double]"6.1225632"
$ws.Range("F16").Value = [double]"1.1431888e-05"
$ws.Range("D17").Value = [double]"53.885328"
$ws.Range("E17").Value = [double]"798.4095600000001"
$ws.Range("F17").Value = [double]"0.0015024659"
$ws.Range("D18").Value = [double]"7.092528"
$ws.Range("E18").Value = [double]"21.110651"
$ws.Range("F18").Value = [double]"0.0001977585"
$ws.Range("D19").Value = [double]"24.01396866666667"
$ws.Range("E19").Value = [double]"41.557433"
$ws.Range("F19").Value = [double]"0.00066957317"
$ws.Range("D20").Value = [double]"6.844884000000001"
$ws.Range("E20").Value = [double]"21.579965"
$ws.Range("F20").Value = [double]"0.00019085352"
$ws.Range("D21").Value = [double]"60.76753466666666"
$ws.Range("E21").Value = [double]"856.94618"
$ws.Range("F21").Value = [double]"0.0016943601"
$ws.Range("D22").Value = [double]"1.4067918"
$ws.Range("E22").Value = [double]"12.031736"
$ws.Range("F22").Value = [double]"3.9225088e-05"
$ws.Range("D23").Value = [double]"1.3067918"
$ws.Range("E23").Value = [double]"12.031736"
$ws.Range("F23").Value = [double]"3.6436822e-05"
$ws.Range("D24").Value = [double]"1.415856533333334"
$ws.Range("E24").Value = [double]"12.73706"
$ws.Range("F24").Value = [double]"3.9477836e-05"
$ws.Range("D25").Value = [double]"1.593752133333334"
$ws.Range("E25").Value = [double]"14.324038"
$ws.Range("F25").Value = [double]"4.4438038e-05"
$ws.Range("D26").Value = [double]"2.133023533333333"
$ws.Range("E26").Value = [double]"18.644144"
$ws.Range("F26").Value = [double]"5.9474355e-05"
$ws.Range("D27").Value = [double]"2.190685266666667"
$ws.Range("E27").Value = [double]"20.407453"
$ws.Range("F27").Value = [double]"6.108211899999999e-05"
$ws.Range("D28").Value = [double]"2.406710333333333"
$ws.Range("E28").Value = [double]"23.405078"
$ws.Range("F28").Value = [double]"6.710546999999999e-05"
$ws.Range("D29").Value = [double]"2.5927354"
$ws.Range("E29").Value = [double]"26.402703"
$ws.Range("F29").Value = [double]"7.229234199999999e-05"
$ws.Range("D30").Value = [double]"10.02228266666667"
$ws.Range("E30").Value = [double]"131.82671"
$ws.Range("F30").Value = [double]"0.00027944783"
$ws.Range("D31").Value = [double]"30.71194"
$ws.Range("E31").Value = [double]"138.8883"
$ws.Range("F31").Value = [double]"0.00085633037"
$ws.Range("D32").Value = [double]"64.788954"
$ws.Range("E32").Value = [double]"229.19117"
$ws.Range("F32").Value = [double]"0.0018064879"
$ws.Range("D33").Value = [double]"5.3832762"
$ws.Range("E33").Value = [double]"44.399003"
$ws.Range("F33").Value = [double]"0.00015010002"
$ws.Range("D34").Value = [double]"3.832836066666667"
$ws.Range("E34").Value = [double]"42.49613"
$ws.Range("F34").Value = [double]"0.00010686964"
$ws.Range("D35").Value = [double]"42.26414933333334"
$ws.Range("E35").Value = [double]"251.23267"
$ws.Range("F35").Value = [double]"0.0011784366"
$ws.Range("D36").Value = [double]"0.1684665933333333"
$ws.Range("E36").Value = [double]"3.0417077"
$ws.Range("F36").Value = [double]"4.6972956e-06"
$ws.Range("D37").Value = [double]"3.418162733333334"
$ws.Range("E37").Value = [double]"61.715809"
$ws.Range("F37").Value = [double]"9.530744599999999e-05"

# --- Add header comments (A1:G1) ---
$ws.Range("A1").AddComment("Data type: Categorical (text)")
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)")
$ws.Range("C1").AddComment("Data type: Categorical (text)")
$ws.Range("D1").AddComment("Data type: Carbon footprint")
$ws.Range("E1").AddComment("Data type: Cumulative energy demand")
$ws.Range("F1").AddComment("Data type: Climate change impact")
$ws.Range("G1").AddComment("Data type: Categorical (text)")
